# Automatic update of files.
# Updates the "Förändrad" (C) date-serial for every data row to 46073,
# and refreshes the latest entries (rows 7-12) with the newest rotation
# of case data in columns A (Beteckning), B (Datum) and G (Area).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 46073

# Row 3
$ws.Range("C3").Value = 46073

# Row 4
$ws.Range("C4").Value = 46073

# Row 5
$ws.Range("C5").Value = 46073

# Row 6
$ws.Range("C6").Value = 46073

# Row 7
$ws.Range("A7").Value = "A 19922-2025"
$ws.Range("B7").Value = 45771.63034722222
$ws.Range("C7").Value = 46073
$ws.Range("G7").Value = 10.1

# Row 8
$ws.Range("A8").Value = "A 25015-2023"
$ws.Range("B8").Value = 45085.6989699074
$ws.Range("C8").Value = 46073
$ws.Range("G8").Value = 1.8

# Row 9
$ws.Range("A9").Value = "A 25634-2025"
$ws.Range("B9").Value = 45803.59570601852
$ws.Range("C9").Value = 46073
$ws.Range("G9").Value = 6

# Row 10
$ws.Range("A10").Value = "A 28266-2025"
$ws.Range("B10").Value = 45818.56381944445
$ws.Range("C10").Value = 46073
$ws.Range("G10").Value = 1.9

# Row 11
$ws.Range("A11").Value = "A 60024-2025"
$ws.Range("B11").Value = 45992
$ws.Range("C11").Value = 46073
$ws.Range("G11").Value = 1.1

# Row 12
$ws.Range("A12").Value = "A 62884-2021"
$ws.Range("B12").Value = 44504
$ws.Range("C12").Value = 46073
$ws.Range("G12").Value = 0.8

# Row 13
$ws.Range("C13").Value = 46073

# Row 14
$ws.Range("C14").Value = 46073
